# Assignment 3 Report: "Z is some constant" -> "Z are some constants"
#
# The sentence currently reads:
#   "... of connections are necessary. Where C, Z is some constant , N is
#    the number of nodes and A is the result."
#
# It should read:
#   "... of connections are necessary. Where C, Z are some constants , N is
#    the number of nodes and A is the result."
#
# i.e. "is some constant" (singular, referring only to Z) becomes
# "are some constants" (plural, since both C and Z are constants).

$d = $word.ActiveDocument

$range = $d.Content
$found = $range.Find.Execute(
    " is some constant",  # Find What
    $true,                # MatchCase
    $false,               # MatchWholeWord
    $false,               # MatchWildcards
    $false,               # MatchSoundsLike
    $false,               # MatchAllWordForms
    $true,                # Forward
    1,                     # Wrap (wdFindContinue)
    $false,               # Format
    " are some constants", # ReplaceWith
    2                      # Replace (wdReplaceAll)
)

if (-not $found) {
    throw "Could not find the target phrase ' is some constant' to replace."
}

Write-Host "Replaced ' is some constant' with ' are some constants'. Find result: $found"
